$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 02:41"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5954816
$ws.Range("C4").Value = 39186
$ws.Range("D4").Value = 3252416
$ws.Range("E4").Value = 2520025
$ws.Range("G4").Value = 1261
$ws.Range("H4").Value = 182375

# Row 15 - Argentina
$ws.Range("B15").Value = 359638
$ws.Range("C15").Value = 8771
$ws.Range("E15").Value = 88873
$ws.Range("G15").Value = 197
$ws.Range("H15").Value = 7563

# Row 23 - Alemania
$ws.Range("B23").Value = 237572
$ws.Range("C23").Value = 1455
$ws.Range("E23").Value = 18627

# Row 27 - Canada
$ws.Range("B27").Value = 125969
$ws.Range("C27").Value = 322
$ws.Range("D27").Value = 112050
$ws.Range("E27").Value = 4829

# Row 36 - Panama
$ws.Range("B36").Value = 88381
$ws.Range("C36").Value = 896
$ws.Range("D36").Value = 62759
$ws.Range("E36").Value = 23703
$ws.Range("G36").Value = 13
$ws.Range("H36").Value = 1919

# Row 59 - Argelia
$ws.Range("D59").Value = 29587
$ws.Range("E59").Value = 11185

# Row 60 - Venezuela
$ws.Range("B60").Value = 41158
$ws.Range("C60").Value = 820
$ws.Range("D60").Value = 32015
$ws.Range("E60").Value = 8800
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 343

# Row 114 - Montenegro
$ws.Range("B114").Value = 4444
$ws.Range("C114").Value = 66
$ws.Range("D114").Value = 3480
$ws.Range("E114").Value = 877
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = 87

# Row 165 - Santo Tome y Principe
$ws.Range("D165").Value = 833
$ws.Range("E165").Value = 44

# Row 175 - Comoras
$ws.Range("D175").Value = 399
$ws.Range("E175").Value = 11

# Row 188 - Bermudas
$ws.Range("D188").Value = 151
$ws.Range("E188").Value = 8
